$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column F entirely (it was empty), which shifts G,H,I left to F,G,H
$ws.Columns.Item(6).Delete()

# Update selection to match target (E1 selected instead of G6)
$ws.Range("E1").Select()
